# Auto-generated Excel COM-interop script applying the Valefor_Profits market-data refresh.
# Updates cached price/profit figures (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value2 = 1556.25  # H19
$ws.Cells.Item(19, 9).Value2 = 2539.6  # I19
$ws.Cells.Item(19, 10).Value2 = 1109.2727  # J19
$ws.Cells.Item(19, 11).Value2 = 2539.6  # K19
$ws.Cells.Item(19, 12).Value2 = 1109.2727  # L19
$ws.Cells.Item(19, 13).Value2 = -2364.6  # M19
$ws.Cells.Item(19, 14).Value2 = -1459.2727  # N19
$ws.Cells.Item(40, 8).Value2 = 5906.087  # H40
$ws.Cells.Item(40, 9).Value2 = 8995.385  # I40
$ws.Cells.Item(40, 10).Value2 = 1890  # J40
$ws.Cells.Item(40, 11).Value2 = 8995.385  # K40
$ws.Cells.Item(40, 12).Value2 = 1890  # L40
$ws.Cells.Item(40, 13).Value2 = -8820.385  # M40
$ws.Cells.Item(40, 14).Value2 = -2240  # N40
$ws.Cells.Item(53, 8).Value2 = 124.23529  # H53
$ws.Cells.Item(53, 9).Value2 = 125.75  # I53
$ws.Cells.Item(53, 11).Value2 = 125.75  # K53
$ws.Cells.Item(53, 13).Value2 = 511.25  # M53
$ws.Cells.Item(62, 8).Value2 = 1377.7778  # H62
$ws.Cells.Item(62, 9).Value2 = 1080  # I62
$ws.Cells.Item(62, 10).Value2 = 1750  # J62
$ws.Cells.Item(62, 11).Value2 = 1080  # K62
$ws.Cells.Item(62, 12).Value2 = 1750  # L62
$ws.Cells.Item(62, 13).Value2 = -456  # M62
$ws.Cells.Item(62, 14).Value2 = -2998  # N62
$ws.Cells.Item(65, 8).Value2 = 1377.7778  # H65
$ws.Cells.Item(65, 9).Value2 = 1080  # I65
$ws.Cells.Item(65, 10).Value2 = 1750  # J65
$ws.Cells.Item(65, 11).Value2 = 5400  # K65
$ws.Cells.Item(65, 12).Value2 = 8750  # L65
$ws.Cells.Item(65, 13).Value2 = -2280  # M65
$ws.Cells.Item(65, 14).Value2 = -14990  # N65
$ws.Cells.Item(125, 8).Value2 = 5510.16  # H125
$ws.Cells.Item(125, 9).Value2 = 5775.778  # I125
$ws.Cells.Item(125, 10).Value2 = 5360.75  # J125
$ws.Cells.Item(125, 11).Value2 = 51982.002  # K125
$ws.Cells.Item(125, 12).Value2 = 48246.75  # L125
$ws.Cells.Item(125, 13).Value2 = -49522.002  # M125
$ws.Cells.Item(125, 14).Value2 = -53166.75  # N125
$ws.Cells.Item(132, 8).Value2 = 3475030  # H132
$ws.Cells.Item(132, 9).Value2 = 5052602.5  # I132
$ws.Cells.Item(132, 10).Value2 = 4371.2  # J132
$ws.Cells.Item(132, 11).Value2 = 15157807.5  # K132
$ws.Cells.Item(132, 12).Value2 = 13113.6  # L132
$ws.Cells.Item(132, 13).Value2 = -15155277.5  # M132
$ws.Cells.Item(132, 14).Value2 = -18173.6  # N132
$ws.Cells.Item(137, 8).Value2 = 1196.9375  # H137
$ws.Cells.Item(137, 9).Value2 = 1143.2  # I137
$ws.Cells.Item(137, 11).Value2 = 3429.6  # K137
$ws.Cells.Item(137, 13).Value2 = -879.6000000000004  # M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value2 = 1550.75  # H2
$ws.Cells.Item(2, 9).Value2 = 1067.6666  # I2
$ws.Cells.Item(2, 10).Value2 = 3000  # J2
$ws.Cells.Item(2, 11).Value2 = 1067.6666  # K2
$ws.Cells.Item(2, 12).Value2 = 3000  # L2
$ws.Cells.Item(2, 13).Value2 = -954.6666  # M2
$ws.Cells.Item(2, 14).Value2 = -3226  # N2
$ws.Cells.Item(32, 8).Value2 = 6489.8394  # H32
$ws.Cells.Item(32, 9).Value2 = 2818.6304  # I32
$ws.Cells.Item(32, 11).Value2 = 2818.6304  # K32
$ws.Cells.Item(32, 13).Value2 = -2531.6304  # M32
$ws.Cells.Item(45, 8).Value2 = 1610.3846  # H45
$ws.Cells.Item(45, 10).Value2 = 1726.2858  # J45
$ws.Cells.Item(45, 12).Value2 = 1726.2858  # L45
$ws.Cells.Item(45, 14).Value2 = -2480.2858  # N45
$ws.Cells.Item(116, 8).Value2 = 1550.75  # H116
$ws.Cells.Item(116, 9).Value2 = 1067.6666  # I116
$ws.Cells.Item(116, 10).Value2 = 3000  # J116
$ws.Cells.Item(116, 11).Value2 = 1067.6666  # K116
$ws.Cells.Item(116, 12).Value2 = 3000  # L116
$ws.Cells.Item(116, 13).Value2 = 1226.3334  # M116
$ws.Cells.Item(116, 14).Value2 = -7588  # N116

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value2 = 1550.75  # H3
$ws.Cells.Item(3, 9).Value2 = 1067.6666  # I3
$ws.Cells.Item(3, 10).Value2 = 3000  # J3
$ws.Cells.Item(3, 11).Value2 = 1067.6666  # K3
$ws.Cells.Item(3, 12).Value2 = 3000  # L3
$ws.Cells.Item(3, 13).Value2 = -953.6666  # M3
$ws.Cells.Item(3, 14).Value2 = -3228  # N3
$ws.Cells.Item(105, 8).Value2 = 4547506.5  # H105
$ws.Cells.Item(105, 9).Value2 = 7577344  # I105
$ws.Cells.Item(105, 11).Value2 = 7577344  # K105
$ws.Cells.Item(105, 13).Value2 = -7575597  # M105
$ws.Cells.Item(126, 8).Value2 = 50768.5  # H126
$ws.Cells.Item(126, 10).Value2 = 50768.5  # J126
$ws.Cells.Item(126, 12).Value2 = 50768.5  # L126
$ws.Cells.Item(126, 14).Value2 = -60648.5  # N126

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 2552.64  # H31
$ws.Cells.Item(31, 9).Value2 = 2100.889  # I31
$ws.Cells.Item(31, 10).Value2 = 3714.2856  # J31
$ws.Cells.Item(31, 11).Value2 = 2100.889  # K31
$ws.Cells.Item(31, 12).Value2 = 3714.2856  # L31
$ws.Cells.Item(31, 13).Value2 = -1805.889  # M31
$ws.Cells.Item(31, 14).Value2 = -4304.2856  # N31
$ws.Cells.Item(34, 8).Value2 = 2552.64  # H34
$ws.Cells.Item(34, 9).Value2 = 2100.889  # I34
$ws.Cells.Item(34, 10).Value2 = 3714.2856  # J34
$ws.Cells.Item(34, 11).Value2 = 2100.889  # K34
$ws.Cells.Item(34, 12).Value2 = 3714.2856  # L34
$ws.Cells.Item(34, 13).Value2 = -1898.889  # M34
$ws.Cells.Item(34, 14).Value2 = -4118.2856  # N34
$ws.Cells.Item(110, 8).Value2 = 40702  # H110
$ws.Cells.Item(110, 10).Value2 = 40702  # J110
$ws.Cells.Item(110, 12).Value2 = 40702  # L110
$ws.Cells.Item(110, 14).Value2 = -48882  # N110
$ws.Cells.Item(132, 8).Value2 = 2910.2  # H132
$ws.Cells.Item(132, 9).Value2 = 1178.5  # I132
$ws.Cells.Item(132, 10).Value2 = 4064.6667  # J132
$ws.Cells.Item(132, 11).Value2 = 3535.5  # K132
$ws.Cells.Item(132, 12).Value2 = 12194.0001  # L132
$ws.Cells.Item(132, 13).Value2 = -1005.5  # M132
$ws.Cells.Item(132, 14).Value2 = -17254.0001  # N132
$ws.Cells.Item(134, 8).Value2 = 2796.7273  # H134
$ws.Cells.Item(134, 9).Value2 = 2207.2  # I134
$ws.Cells.Item(134, 10).Value2 = 3288  # J134
$ws.Cells.Item(134, 11).Value2 = 6621.599999999999  # K134
$ws.Cells.Item(134, 12).Value2 = 9864  # L134
$ws.Cells.Item(134, 13).Value2 = -4086.599999999999  # M134
$ws.Cells.Item(134, 14).Value2 = -14934  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value2 = 31.631578  # H12
$ws.Cells.Item(12, 9).Value2 = 19.571428  # I12
$ws.Cells.Item(12, 11).Value2 = 58.71428400000001  # K12
$ws.Cells.Item(12, 13).Value2 = 114.285716  # M12
$ws.Cells.Item(23, 8).Value2 = 222.8  # H23
$ws.Cells.Item(23, 9).Value2 = 231.16667  # I23
$ws.Cells.Item(23, 10).Value2 = 219.21428  # J23
$ws.Cells.Item(23, 11).Value2 = 693.50001  # K23
$ws.Cells.Item(23, 12).Value2 = 657.64284  # L23
$ws.Cells.Item(23, 13).Value2 = -458.50001  # M23
$ws.Cells.Item(23, 14).Value2 = -1127.64284  # N23
$ws.Cells.Item(39, 8).Value2 = 2582.0688  # H39
$ws.Cells.Item(39, 9).Value2 = 745  # I39
$ws.Cells.Item(39, 10).Value2 = 2876  # J39
$ws.Cells.Item(39, 11).Value2 = 2235  # K39
$ws.Cells.Item(39, 12).Value2 = 8628  # L39
$ws.Cells.Item(39, 13).Value2 = -1941  # M39
$ws.Cells.Item(39, 14).Value2 = -9216  # N39
$ws.Cells.Item(68, 8).Value2 = 3679.1428  # H68
$ws.Cells.Item(68, 9).Value2 = 876  # I68
$ws.Cells.Item(68, 10).Value2 = 4800.4  # J68
$ws.Cells.Item(68, 11).Value2 = 2628  # K68
$ws.Cells.Item(68, 12).Value2 = 14401.2  # L68
$ws.Cells.Item(68, 13).Value2 = -1817  # M68
$ws.Cells.Item(68, 14).Value2 = -16023.2  # N68
$ws.Cells.Item(71, 8).Value2 = 3679.1428  # H71
$ws.Cells.Item(71, 9).Value2 = 876  # I71
$ws.Cells.Item(71, 10).Value2 = 4800.4  # J71
$ws.Cells.Item(71, 11).Value2 = 7884  # K71
$ws.Cells.Item(71, 12).Value2 = 43203.6  # L71
$ws.Cells.Item(71, 13).Value2 = -3828  # M71
$ws.Cells.Item(71, 14).Value2 = -51315.6  # N71
$ws.Cells.Item(86, 8).Value2 = 306.66666  # H86
$ws.Cells.Item(86, 9).Value2 = 290  # I86
$ws.Cells.Item(86, 11).Value2 = 870  # K86
$ws.Cells.Item(86, 13).Value2 = 316  # M86
$ws.Cells.Item(89, 8).Value2 = 306.66666  # H89
$ws.Cells.Item(89, 9).Value2 = 290  # I89
$ws.Cells.Item(89, 11).Value2 = 2610  # K89
$ws.Cells.Item(89, 13).Value2 = 3318  # M89
$ws.Cells.Item(109, 8).Value2 = 1099.3334  # H109
$ws.Cells.Item(109, 9).Value2 = 149  # I109
$ws.Cells.Item(109, 10).Value2 = 3000  # J109
$ws.Cells.Item(109, 11).Value2 = 447  # K109
$ws.Cells.Item(109, 12).Value2 = 9000  # L109
$ws.Cells.Item(109, 13).Value2 = 593  # M109
$ws.Cells.Item(109, 14).Value2 = -11080  # N109
$ws.Cells.Item(122, 8).Value2 = 4983.2915  # H122
$ws.Cells.Item(122, 9).Value2 = 590.6667  # I122
$ws.Cells.Item(122, 10).Value2 = 6447.5  # J122
$ws.Cells.Item(122, 11).Value2 = 5316.0003  # K122
$ws.Cells.Item(122, 12).Value2 = 58027.5  # L122
$ws.Cells.Item(122, 13).Value2 = -2866.0003  # M122
$ws.Cells.Item(122, 14).Value2 = -62927.5  # N122
$ws.Cells.Item(131, 8).Value2 = 957.9583  # H131
$ws.Cells.Item(131, 10).Value2 = 1026  # J131
$ws.Cells.Item(131, 12).Value2 = 3078  # L131
$ws.Cells.Item(131, 14).Value2 = -13158  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value2 = 1490.3549  # H97
$ws.Cells.Item(97, 9).Value2 = 1294.5264  # I97
$ws.Cells.Item(97, 10).Value2 = 1800.4166  # J97
$ws.Cells.Item(97, 11).Value2 = 1294.5264  # K97
$ws.Cells.Item(97, 12).Value2 = 1800.4166  # L97
$ws.Cells.Item(97, 13).Value2 = -798.5264  # M97
$ws.Cells.Item(97, 14).Value2 = -2792.4166  # N97
$ws.Cells.Item(102, 8).Value2 = 3445.739  # H102
$ws.Cells.Item(102, 9).Value2 = 2452.7144  # I102
$ws.Cells.Item(102, 11).Value2 = 2452.7144  # K102
$ws.Cells.Item(102, 13).Value2 = -830.7143999999998  # M102

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 2081.25  # H7
$ws.Cells.Item(7, 9).Value2 = 2107.6924  # I7
$ws.Cells.Item(7, 10).Value2 = 1966.6666  # J7
$ws.Cells.Item(7, 11).Value2 = 2107.6924  # K7
$ws.Cells.Item(7, 12).Value2 = 1966.6666  # L7
$ws.Cells.Item(7, 13).Value2 = -1995.6924  # M7
$ws.Cells.Item(7, 14).Value2 = -2190.6666  # N7
$ws.Cells.Item(16, 8).Value2 = 1225.9  # H16
$ws.Cells.Item(16, 9).Value2 = 1225.9  # I16
$ws.Cells.Item(16, 10).Value2 = 0  # J16
$ws.Cells.Item(16, 11).Value2 = 1225.9  # K16
$ws.Cells.Item(16, 12).Value2 = 0  # L16
$ws.Cells.Item(16, 13).Value2 = -1055.9  # M16
$ws.Cells.Item(16, 14).ClearContents()  # N16
$ws.Cells.Item(22, 8).Value2 = 423.4516  # H22
$ws.Cells.Item(22, 9).Value2 = 411.67856  # I22
$ws.Cells.Item(22, 10).Value2 = 533.3333  # J22
$ws.Cells.Item(22, 11).Value2 = 411.67856  # K22
$ws.Cells.Item(22, 12).Value2 = 533.3333  # L22
$ws.Cells.Item(22, 13).Value2 = -116.67856  # M22
$ws.Cells.Item(22, 14).Value2 = -1123.3333  # N22
$ws.Cells.Item(27, 8).Value2 = 423.4516  # H27
$ws.Cells.Item(27, 9).Value2 = 411.67856  # I27
$ws.Cells.Item(27, 10).Value2 = 533.3333  # J27
$ws.Cells.Item(27, 11).Value2 = 411.67856  # K27
$ws.Cells.Item(27, 12).Value2 = 533.3333  # L27
$ws.Cells.Item(27, 13).Value2 = -304.67856  # M27
$ws.Cells.Item(27, 14).Value2 = -747.3333  # N27
$ws.Cells.Item(126, 8).Value2 = 2081.25  # H126
$ws.Cells.Item(126, 9).Value2 = 2107.6924  # I126
$ws.Cells.Item(126, 10).Value2 = 1966.6666  # J126
$ws.Cells.Item(126, 11).Value2 = 6323.0772  # K126
$ws.Cells.Item(126, 12).Value2 = 5899.9998  # L126
$ws.Cells.Item(126, 13).Value2 = -3853.0772  # M126
